# Spacecraft Operational Modes.xlsx - update power draw figures and
# rename the "Sintering compression" activity to a generic "Compression"
# activity that now occurs during Regolith collection instead of Sintering.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mode Descriptions")

# --- Furnace row (row 5): Sintering power draw increases 100 -> 750 ---
$ws.Range("H5").Value = 750

# --- Compressing system row (row 6) ---
# The compression activity moves from the Sintering columns (G/H) to the
# Regolith collection columns (E/F), and is renamed "Compression".
$ws.Range("E6").Value = "Compression"
$ws.Range("F6").Value = 140
$ws.Range("G6").Value = "off"
$ws.Range("H6").Value = 0

# Testing power draw for compressing system increases 100 -> 140
$ws.Range("L6").Value = 140
